$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sending / target cluster labels introduced by this edit.
# Ligand symbol (B) stays "Cxcl5", Receptor symbol (C) stays "Cxcr2" for all rows.
# Sending cluster (A) and Target cluster (D) vary across the new rows.

$data = @(
    @{ A="FAPs"; D="ECs";    E=3; F=1;                  G=319.3801833333334;  H=958.1405500000001; I=0.9980732501448392;   J=0.9980732501448393;   K=2; L=0.6666666666666666; M=0.040154;            N=0.120462;           O=0.000117739357997762;  P=0.000117739357997762;  Q=12.82439188156667;   R=115.4195269341;      S=[double]"0.0001175125037067931";  T=[double]"0.0001175125037067931" },
    @{ A="FAPs"; D="M1";     E=3; F=1;                  G=319.3801833333334;  H=958.1405500000001; I=0.9980732501448392;   J=0.9980732501448393;   K=2; L=0.6666666666666666; M=2.854025333333333;   N=8.562075999999999;  O=0.008368558809981951;  P=0.008368558809981953;  Q=911.5191341979778;   R=8203.672207781799;   S=[double]"0.008352434690506914";   T=[double]"0.008352434690506916" },
    @{ A="FAPs"; D="M2";     E=3; F=1;                  G=319.3801833333334;  H=958.1405500000001; I=0.9980732501448392;   J=0.9980732501448393;   K=2; L=0.6666666666666666; M=0.1353773333333333;  N=0.406132;           O=0.0003969527398046444; P=0.0003969527398046444; Q=43.23683753917778;   R=389.1315378526;      S=[double]"0.0003961879111707201";  T=[double]"0.0003961879111707201" },
    @{ A="FAPs"; D="Neutro"; E=3; F=1;                  G=319.3801833333334;  H=958.1405500000001; I=0.9980732501448392;   J=0.9980732501448393;   K=3; L=1;                  M=338.0118816666666;   N=1014.035645;        O=0.9911167490922156;    P=0.9911167490922157;    Q=107954.296735545;    R=971588.6706199048;   S=[double]"0.9892071150394547";     T=[double]"0.9892071150394549" },
    @{ A="sCs";  D="ECs";    E=1; F=0.3333333333333333; G=0.6165536666666667; H=1.849661;          I=0.001926749855160762; J=0.001926749855160763; K=2; L=0.6666666666666666; M=0.040154;            N=0.120462;           O=0.000117739357997762;  P=0.000117739357997762;  Q=0.02475709593133333; R=0.222813863382;      S=[double]"2.26854290968909E-07";   T=[double]"2.268542909689091E-07" },
    @{ A="sCs";  D="M1";     E=1; F=0.3333333333333333; G=0.6165536666666667; H=1.849661;          I=0.001926749855160762; J=0.001926749855160763; K=2; L=0.6666666666666666; M=2.854025333333333;   N=8.562075999999999;  O=0.008368558809981951;  P=0.008368558809981953;  Q=1.759659784026222;   R=15.836938056236;     S=[double]"1.612411947503704E-05";  T=[double]"1.612411947503705E-05" },
    @{ A="sCs";  D="M2";     E=1; F=0.3333333333333333; G=0.6165536666666667; H=1.849661;          I=0.001926749855160762; J=0.001926749855160763; K=2; L=0.6666666666666666; M=0.1353773333333333;  N=0.406132;           O=0.0003969527398046444; P=0.0003969527398046444; Q=0.08346739125022222; R=0.7512065212519999;  S=[double]"7.648286339242662E-07";  T=[double]"7.648286339242665E-07" },
    @{ A="sCs";  D="Neutro"; E=1; F=0.3333333333333333; G=0.6165536666666667; H=1.849661;          I=0.001926749855160762; J=0.001926749855160763; K=3; L=1;                  M=338.0118816666666;   N=1014.035645;        O=0.9911167490922156;    P=0.9911167490922157;    Q=208.4024650184828;   R=1875.622185166345;   S=[double]"0.001909634052760832";   T=[double]"0.001909634052760832" }
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec.A
    $ws.Cells.Item($rowIndex, 2).Value = "Cxcl5"
    $ws.Cells.Item($rowIndex, 3).Value = "Cxcr2"
    $ws.Cells.Item($rowIndex, 4).Value = $rec.D
    $ws.Cells.Item($rowIndex, 5).Value = $rec.E
    $ws.Cells.Item($rowIndex, 6).Value = $rec.F
    $ws.Cells.Item($rowIndex, 7).Value = $rec.G
    $ws.Cells.Item($rowIndex, 8).Value = $rec.H
    $ws.Cells.Item($rowIndex, 9).Value = $rec.I
    $ws.Cells.Item($rowIndex, 10).Value = $rec.J
    $ws.Cells.Item($rowIndex, 11).Value = $rec.K
    $ws.Cells.Item($rowIndex, 12).Value = $rec.L
    $ws.Cells.Item($rowIndex, 13).Value = $rec.M
    $ws.Cells.Item($rowIndex, 14).Value = $rec.N
    $ws.Cells.Item($rowIndex, 15).Value = $rec.O
    $ws.Cells.Item($rowIndex, 16).Value = $rec.P
    $ws.Cells.Item($rowIndex, 17).Value = $rec.Q
    $ws.Cells.Item($rowIndex, 18).Value = $rec.R
    $ws.Cells.Item($rowIndex, 19).Value = $rec.S
    $ws.Cells.Item($rowIndex, 20).Value = $rec.T
    $rowIndex++
}
